$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 2.31998719698953
$ws.Range("C2").Value = 0.3284434136486425
$ws.Range("E2").Value = 0.06291430085788541
$ws.Range("F2").Value = 3.344146860200027
$ws.Range("G2").Value = 0.002568520650615289
$ws.Range("I2").Value = 2.070593412977615
$ws.Range("J2").Value = 0.1552357495912808
$ws.Range("L2").Value = 0.3815221499667274

# Row 3
$ws.Range("B3").Value = 2.22054201590953
$ws.Range("C3").Value = 0.293117759577683
$ws.Range("E3").Value = 0.06249437285523918
$ws.Range("F3").Value = 3.308271484791049
$ws.Range("G3").Value = 0.002574637840523529
$ws.Range("I3").Value = 2.052065791112952
$ws.Range("J3").Value = 0.1555323562710349
$ws.Range("L3").Value = 0.3778741111961637

# Row 4
$ws.Range("B4").Value = 2.160984004844067
$ws.Range("C4").Value = 0.2715803440934224
$ws.Range("E4").Value = 0.06223419967362975
$ws.Range("F4").Value = 3.287932731733747
$ws.Range("G4").Value = 0.002578588976247293
$ws.Range("I4").Value = 2.04169652986873
$ws.Range("J4").Value = 0.1557677312900161
$ws.Range("L4").Value = 0.3758282970611191

# Row 5
$ws.Range("B5").Value = 2.137090391879667
$ws.Range("C5").Value = 0.2628410106868557
$ws.Range("E5").Value = 0.06212755642720857
$ws.Range("F5").Value = 3.280067267120515
$ws.Range("G5").Value = 0.002580248345710868
$ws.Range("I5").Value = 2.037722450912369
$ws.Range("J5").Value = 0.1558769292900308
$ws.Range("L5").Value = 0.3750435677590218

# Row 6
$ws.Range("B6").Value = 2.133145607069764
$ws.Range("C6").Value = 0.2613920756346886
$ws.Range("E6").Value = 0.06210980995323645
$ws.Range("F6").Value = 3.278786688962171
$ws.Range("G6").Value = 0.002580526862508141
$ws.Range("I6").Value = 2.037077697783786
$ws.Range("J6").Value = 0.1558958604131391
$ws.Range("L6").Value = 0.3749162253863219

# Row 7
$ws.Range("B7").Value = 2.160660242822644
$ws.Range("C7").Value = 0.2714623325886123
$ws.Range("E7").Value = 0.06223276400181277
$ws.Range("F7").Value = 3.287824946197119
$ws.Range("G7").Value = 0.002578611155496956
$ws.Range("I7").Value = 2.041641918031885
$ws.Range("J7").Value = 0.155769150345705
$ws.Range("L7").Value = 0.3758175155011827

# Row 8
$ws.Range("B8").Value = 2.28538624032177
$ws.Range("C8").Value = 0.316230791001658
$ws.Range("E8").Value = 0.06276997425497477
$ws.Range("F8").Value = 3.331425316625555
$ws.Range("G8").Value = 0.002570589459120823
$ws.Range("I8").Value = 2.063995109861324
$ws.Range("J8").Value = 0.1553268988948773
$ws.Range("L8").Value = 0.3802241077392807

# Row 9
$ws.Range("B9").Value = 2.541946728093819
$ws.Range("C9").Value = 0.4052853855436069
$ws.Range("E9").Value = 0.06380645461354018
$ws.Range("F9").Value = 3.430421602009119
$ws.Range("G9").Value = 0.002556399365506066
$ws.Range("I9").Value = 2.115897865763003
$ws.Range("J9").Value = 0.1548872783843649
$ws.Range("L9").Value = 0.3904000325702839

# Row 10
$ws.Range("B10").Value = 2.737847745461352
$ws.Range("C10").Value = 0.4715636487910615
$ws.Range("E10").Value = 0.06455980792105809
$ws.Range("F10").Value = 3.511526093999436
$ws.Range("G10").Value = 0.002546901692801917
$ws.Range("I10").Value = 2.159065937826952
$ws.Range("J10").Value = 0.1548320575915696
$ws.Range("L10").Value = 0.3988061308000965

# Row 11
$ws.Range("B11").Value = 2.828602061844606
$ws.Range("C11").Value = 0.5019176601953745
$ws.Range("E11").Value = 0.06490124586905566
$ws.Range("F11").Value = 3.550272692147928
$ws.Range("G11").Value = 0.002542779998705837
$ws.Range("I11").Value = 2.179823177440554
$ws.Range("J11").Value = 0.1548666194728057
$ws.Range("L11").Value = 0.4028311602997832

# Row 12
$ws.Range("B12").Value = 2.86320558592007
$ws.Range("C12").Value = 0.5134425565116203
$ws.Range("E12").Value = 0.06503039683906398
$ws.Range("F12").Value = 3.565213724725965
$ws.Range("G12").Value = 0.002541247628261095
$ws.Range("I12").Value = 2.187846427649575
$ws.Range("J12").Value = 0.1548884154478642
$ws.Range("L12").Value = 0.4043841382013085

# Row 13
$ws.Range("B13").Value = 2.855742543187318
$ws.Range("C13").Value = 0.5109590884719069
$ws.Range("E13").Value = 0.06500258766861577
$ws.Range("F13").Value = 3.561983926538971
$ws.Range("G13").Value = 0.002541576390095282
$ws.Range("I13").Value = 2.186111200207492
$ws.Range("J13").Value = 0.154883331980038
$ws.Range("L13").Value = 0.4040483982685856

# Row 14
$ws.Range("B14").Value = 2.831444160402782
$ws.Range("C14").Value = 0.502865202321459
$ws.Range("E14").Value = 0.06491187386684771
$ws.Range("F14").Value = 3.551496504988478
$ws.Range("G14").Value = 0.002542653360877667
$ws.Range("I14").Value = 2.180479978689391
$ws.Range("J14").Value = 0.1548682373915611
$ws.Range("L14").Value = 0.4029583484977763

# Row 15
$ws.Range("B15").Value = 2.816591581858916
$ws.Range("C15").Value = 0.4979114701276899
$ws.Range("E15").Value = 0.06485629140517446
$ws.Range("F15").Value = 3.54510769533681
$ws.Range("G15").Value = 0.002543316733322126
$ws.Range("I15").Value = 2.177051968365376
$ws.Range("J15").Value = 0.1548601293335352
$ws.Range("L15").Value = 0.4022944067822039

# Row 16
$ws.Range("B16").Value = 2.731949806681257
$ws.Range("C16").Value = 0.469584143742793
$ws.Range("E16").Value = 0.06453747189441827
$ws.Range("F16").Value = 3.509031367003303
$ws.Range("G16").Value = 0.002547175042338389
$ws.Range("I16").Value = 2.157732099901395
$ws.Range("J16").Value = 0.1548310105867898
$ws.Range("L16").Value = 0.39854712060297

# Row 17
$ws.Range("B17").Value = 2.680445240744575
$ws.Range("C17").Value = 0.4522592653268021
$ws.Range("E17").Value = 0.06434158890835384
$ws.Range("F17").Value = 3.487375583992247
$ws.Range("G17").Value = 0.00254959279917507
$ws.Range("I17").Value = 2.146168156781343
$ws.Range("J17").Value = 0.1548285224453423
$ws.Range("L17").Value = 0.3962996863282484

# Row 18
$ws.Range("B18").Value = 2.650975267344222
$ws.Range("C18").Value = 0.4423135534465814
$ws.Range("E18").Value = 0.06422880070339954
$ws.Range("F18").Value = 3.475093910023162
$ws.Range("G18").Value = 0.002551002155221449
$ws.Range("I18").Value = 2.139622175278078
$ws.Range("J18").Value = 0.1548327013006627
$ws.Range("L18").Value = 0.3950259528233317

# Row 19
$ws.Range("B19").Value = 2.64102366756174
$ws.Range("C19").Value = 0.4389493505170208
$ws.Range("E19").Value = 0.06419059054477749
$ws.Range("F19").Value = 3.470965390027231
$ws.Range("G19").Value = 0.002551482560277392
$ws.Range("I19").Value = 2.137423844948955
$ws.Range("J19").Value = 0.1548350758872274
$ws.Range("L19").Value = 0.3945979438772582

# Row 20
$ws.Range("B20").Value = 2.685912032745989
$ws.Range("C20").Value = 0.4541015428119977
$ws.Range("E20").Value = 0.06436245328789703
$ws.Range("F20").Value = 3.489662837954256
$ws.Range("G20").Value = 0.0025493334879411
$ws.Range("I20").Value = 2.147388247174604
$ws.Range("J20").Value = 0.154828205793514
$ws.Range("L20").Value = 0.3965369711489046

# Row 21
$ws.Range("B21").Value = 2.838574746124095
$ws.Range("C21").Value = 0.5052417357046011
$ws.Range("E21").Value = 0.06493852231675312
$ws.Range("F21").Value = 3.554569608778763
$ws.Range("G21").Value = 0.00254233625832265
$ws.Range("I21").Value = 2.182129567752696
$ws.Range("J21").Value = 0.1548724336923755
$ws.Range("L21").Value = 0.4032777422494291

# Row 22
$ws.Range("B22").Value = 2.939729757987834
$ws.Range("C22").Value = 0.5388431042305797
$ws.Range("E22").Value = 0.0653141944561515
$ws.Range("F22").Value = 3.598556241405959
$ws.Range("G22").Value = 0.002537928775260452
$ws.Range("I22").Value = 2.205785626726538
$ws.Range("J22").Value = 0.1549521561381013
$ws.Range("L22").Value = 0.4078509615335975

# Row 23
$ws.Range("B23").Value = 2.885614606026024
$ws.Range("C23").Value = 0.5208927056391985
$ws.Range("E23").Value = 0.06511375370617145
$ws.Range("F23").Value = 3.574935677541674
$ws.Range("G23").Value = 0.002540266033708456
$ws.Range("I23").Value = 2.193072331342734
$ws.Range("J23").Value = 0.154904915109995
$ws.Range("L23").Value = 0.4053948404811223

# Row 24
$ws.Range("B24").Value = 2.683440057967971
$ws.Range("C24").Value = 0.4532686038025986
$ws.Range("E24").Value = 0.06435302104205398
$ws.Range("F24").Value = 3.488628245693008
$ws.Range("G24").Value = 0.002549450662153402
$ws.Range("I24").Value = 2.146836325861045
$ws.Range("J24").Value = 0.1548283314891279
$ws.Range("L24").Value = 0.3964296375270209

# Row 25
$ws.Range("B25").Value = 2.471247711878448
$ws.Range("C25").Value = 0.3810503207228066
$ws.Range("E25").Value = 0.06352765182307518
$ws.Range("F25").Value = 3.402180806766978
$ws.Range("G25").Value = 0.002560074418287771
$ws.Range("I25").Value = 2.100981044577196
$ws.Range("J25").Value = 0.1549596987424202
$ws.Range("L25").Value = 0.3874834901034205
